$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the cell formats of B1 and D1 (D1 currently carries the red
# "invalid header" style) before changing their text, so the red
# style follows the "INVALID" text to its new location (B1), and D1
# reverts back to the plain default style.
$ws.Range("D1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B1").Value = "INVALID"
$ws.Range("D1").Value = "Target Object"

$ws.Range("D13").Select() | Out-Null
